$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to make room for "SiteCode"
$ws.Columns("B:B").Insert()
$ws.Columns("B:B").ColumnWidth = 11.166666666666666

# Update header row: A1 becomes "Site", B1 becomes "SiteCode"
$ws.Range("A1").Value = "Site"
$ws.Range("B1").Value = "SiteCode"
$ws.Rows("1:1").RowHeight = 32

# Fill in the new SiteCode column (B2:B12) to match each site in column A
# (order chosen to mirror the original authoring sequence of site codes)
$ws.Range("B2").Value = "BBWM"   # Bear Brook
$ws.Range("B4").Value = "DOR"    # Dorset
$ws.Range("B6").Value = "HBEF"   # Hubbard Brook
$ws.Range("B7").Value = "HJA"    # HJ Andrews
$ws.Range("B9").Value = "MEF"    # Marcell
$ws.Range("B11").Value = "SLP"   # Sleepers
$ws.Range("B12").Value = "TLW"   # Turkey Lakes
$ws.Range("B3").Value = "COW"    # Coweeta
$ws.Range("B8").Value = "LUQ"    # Luquilllo
$ws.Range("B10").Value = "SAN"   # Santee
$ws.Range("B5").Value = "ELA"    # ELA

# Restore the active selection to the cell the original author left selected
[void]$ws.Range("B11").Select()
